$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New content for rows 2-16 (job_name, slurm_file_name, result_file_name)
$data = @(
    @("ppta_ptmcmc_pe_gwb_fixsl_infer_orf_set_x_1_ephem_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_pe_gwb_fixslope_infer_orf_set_x_1_ephem_0_20201022.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20201020_gwb_fixgamma_varorf_x1psr/completed.txt"),
    @("ppta_ptmcmc_pe_gwb_cpl_fixslope_set_x_1_eph_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_pe_cpl_gwb_fixslope_set_x_1_ephem_0_20201103.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20200908_cpl_gwb_fixgamma_x1psr/completed.txt"),
    @("ppta_ptmcmc_pe_gwb_cpl_fixsl_inforf_set_x_1_eph_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_pe_gwb_cpl_fixslope_infer_orf_set_x_1_ephem_0_20201030.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20201030_gwb_cpl_fixgamma_varorf_x1psr/completed.txt"),
    @("ppta_ptmcmc_ms_gwb_cpl_fixslope_set_x_1_ephem_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_ms_cpl_gwb_fixslope_set_x_1_ephem_0_20201027.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20200910_cpl_fixgam_20200908_gwb_fixgamma_x1psr/completed.txt"),
    @("ppta_ptmcmc_ms_cpl_gwb_cpl_fixslope_set_x_1_eph_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_ms_cpl_cpl_gwb_fixslope_set_x_1_ephem_0_20201030.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20201030_gwb_fixgamma_namehd_20200908_cpl_gwb_fixgamma_x1psr/completed.txt"),
    @("ppta_ptmcmc_pe_gwb_fixsl_inter_orf_set_x_1_ephem_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_pe_gwb_fixslope_interp_orf_set_x_1_ephem_0_20201105.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20201020_gwb_fixgamma_interporf_x1psr/completed.txt"),
    @("ppta_ptmcmc_pe_gwb_cpl_fixsl_intorf_set_x_1_eph_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_pe_gwb_cpl_fixslope_interp_orf_set_x_1_ephem_0_20201105.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20201105_gwb_cpl_fixgamma_interporf_x1psr/completed.txt"),
    @("ppta_ptmcmc_pe_gwb_cpl_noa_fixsl_set_x_1_eph_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_pe_cpl_gwb_noauto_fixslope_set_20_nf_x_1_ephem_0_20201109.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20200908_cpl_gwb_noa_fixgam_20_nf_x1psr/completed.txt"),
    @("ppta_ptmcmc_pe_gwb_fixsl_inforf_20nf_set_x1_eph_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_pe_gwb_fixslope_infer_orf_20_nf_set_x_1_ephem_0_20201111.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20201020_gwb_fixgamma_varorf_20_nf_x1psr/completed.txt"),
    @("ppta_ptmcmc_ms_gwb_cpl_fixsl_20nf_set_x1_eph0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_ms_cpl_gwb_fixslope_20_nf_set_x_1_ephem_0_20201111.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20200910_cpl_fixgam_20_nf_20200908_gwb_fixgamma_20_nf_x1psr/completed.txt"),
    @("ppta_ptmcmc_pe_gwb_fixsl_intorf_20nf_setx1_eph0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_pe_gwb_fixslope_interp_orf_20_nf_set_x_1_ephem_0_20201105.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20201111_gwb_fixgamma_interporf_20_nf_x1psr/completed.txt"),
    @("ppta_ptmcmc_ms_mon_cpl_fixsl_20nf_set_x_1", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_ms_cpl_monop_fixslope_20_nf_set_x_1_ephem_0_20201117.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20200910_cpl_fixgam_20_nf_20201117_monop_fixgamma_20_nf_x1psr/completed.txt"),
    @("ppta_ptmcmc_ms_dip_cpl_fixsl_20nf_set_x_1", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_ms_cpl_dipol_fixslope_20_nf_set_x_1_ephem_0_20201117.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20200910_cpl_fixgam_20_nf_20201117_dipol_fixgamma_20_nf_x1psr/completed.txt"),
    @("ppta_ptmcmc_ms_dip_cpl_fixslope_set_3_1_ephem_0", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_ms_cpl_dipol_fixslope_set_3_1_ephem_0_20201121.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20200910_cpl_fixgam_20201103_dipol_fixgamma_3psr/completed.txt"),
    @("ppta_ptmcmc_ms_hdip_cpl_fixsl_20nf_set_x1", "/home/bgonchar/correlated_noise_pta_2020/slurm/ppta_dr2_ptmcmc_ms_cpl_halfdip_fixslope_20_nf_set_x_1_20201121.sh", "/fred/oz002/bgoncharov/correlated_noise_pta_2020_out/dr2_timing_20200607/20200910_cpl_fixgam_20_nf_20201121_halfdip_fixgamma_20_nf_x1psr/completed.txt")
)

# Write the new values into rows 2-16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Delete the old rows 17-24 (now leftover beyond the new data range)
$ws.Rows("17:24").Delete()
